$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 777, shifting the existing rows 777:819 down to 778:820.
$ws.Rows.Item(777).Insert()

# Populate the newly inserted row 777 with the new price-report record.
$ws.Range("A777").Value = 4
$ws.Range("B777").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C777").Value = "Los Lagos"
$ws.Range("D777").Value = 45267
$ws.Range("E777").Value = 10
$ws.Range("F777").Value = 100114001
$ws.Range("G777").Value = "Papa"
$ws.Range("H777").Value = "Patagonia"
$ws.Range("I777").Value = "1a nueva(o)"
$ws.Range("J777").Value = 500
$ws.Range("K777").Value = 25000
$ws.Range("L777").Value = 27000
$ws.Range("M777").Value = 26000
$ws.Range("N777").Value = "$/saco 25 kilos"
$ws.Range("O777").Value = "Región de La Araucanía"
$ws.Range("P777").Value = 1040
$ws.Range("Q777").Value = 25
$ws.Range("R777").Value = "Hortaliza"
